# "harmonized similar tags to be the same"
#
# The #TAGS list on the "SwateTemplateMetadata" sheet mixed free-text
# tags (some with stray leading/trailing spaces, one a bare ontology URL)
# with the rest of the tag set. This normalizes them:
#   - drops the redundant " protocol " tag entirely
#   - trims the stray whitespace from "phenotyping " / "study "
#   - replaces the raw NCIT URL + "NCIT" pair for the "Plant" tag with a
#     short CURIE ("NCIT:C14258")
#   - gives the "phenotyping" and "study" tags their own term accession
#     numbers ("DPBO:1000224" / "NCIT:C63536") which they previously
#     lacked

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("SwateTemplateMetadata")

# --- Row 13/14: "Tags Term Accession Number" / "Tags Term Source REF".
#     Carry the wrap-text formatting that lived on E13/E14 over to the
#     new D13/D14 position (still the "Plant" column) before touching
#     any values.
$ws.Range("E13").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null

# "Plant" tag's accession: the old raw NCIT URL becomes a short CURIE.
$ws.Range("D13").Value = "NCIT:C14258"

# --- Row 12: "#TAGS list" -> "Tags" values, shifted one column left
#     after dropping the old " protocol " tag, with "phenotyping "/
#     "study " trimmed of their stray whitespace.
$ws.Range("C12").Value = "phenotyping"
$ws.Range("D12").Value = "Plant"
$ws.Range("E12").Value = " metadata "

# "phenotyping" tag gets its own (previously missing) term accession.
$ws.Range("C13").Value = "DPBO:1000224"

$ws.Range("F12").Value = "study"
$ws.Range("G12").Value = "MIAPPE"
$ws.Range("H12").ClearContents()

# "study" tag gets its own (previously missing) term accession.
$ws.Range("F13").Value = "NCIT:C63536"

# "Plant" tag's term source REF is folded into the CURIE above now, so
# its own column is left blank (still wrap-text formatted).
$ws.Range("D14").ClearContents()

# Old "Plant" accession/source-ref cells (now vacated) are cleared.
$ws.Range("E13").Clear() | Out-Null
$ws.Range("E14").Clear() | Out-Null

# Row 13 only needs two stacked lines now (vs. the long URL before),
# so it no longer needs to be as tall.
$ws.Rows.Item(13).RowHeight = 28.8

$ws.Activate() | Out-Null
$ws.Range("G18").Select() | Out-Null

$excel.CalculateFull() | Out-Null
